# "change state to simple Column"
#
# The original "Status" sheet (3 columns: migration / stateId / stateName)
# is demoted to a trailing sheet named "Status_1" (its data is untouched),
# and a brand-new "Status" sheet takes its old place with a simplified,
# two-column layout (migration / stateName only - the numeric stateId
# column is dropped). "Permissions" also gains a new lookup row and
# becomes the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new sheet immediately before the existing "Status"
#    sheet. Excel hands it the next free sheetId and slots it in at
#    position 1; the original "Status" sheet (old 3-column data) is
#    pushed down to position 2.
# ---------------------------------------------------------------------
$oldStatusRef = $wb.Worksheets.Item("Status")
$null = $wb.Worksheets.Add($oldStatusRef)

# References can shift after structural edits (insert/move/rename), so
# re-fetch both sheets positionally right after the Add.
$newStatusSheet = $wb.Worksheets.Item(1)
$oldStatusSheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 2. Free up the "Status" name: rename the old sheet to "Status_1" first,
#    then claim "Status" for the freshly inserted sheet.
# ---------------------------------------------------------------------
$oldStatusSheet.Name = "Status_1"
$newStatusSheet.Name = "Status"

# ---------------------------------------------------------------------
# 3. Populate the new "Status" sheet with the simplified two-column
#    layout (migration flag + plain state name, no numeric stateId).
# ---------------------------------------------------------------------
$newStatusSheet.Range("A1").Value = "migration"
$newStatusSheet.Range("B1").Value = "stateName"
$newStatusSheet.Range("A2").Value = $true
$newStatusSheet.Range("B2").Value = "Active"
$newStatusSheet.Range("A3").Value = $false
$newStatusSheet.Range("B3").Value = "Deactivate"
$newStatusSheet.Range("A4").Value = $false
$newStatusSheet.Range("B4").Value = "Delete"

# Carry over the header/text styling (style index 1) from the old sheet
# onto the matching cells of the new layout.
$oldStatusSheet.Range("A1").Copy()
$newStatusSheet.Range("A1").PasteSpecial(-4122)
$oldStatusSheet.Range("A1").Copy()
$newStatusSheet.Range("B1").PasteSpecial(-4122)
$oldStatusSheet.Range("C2").Copy()
$newStatusSheet.Range("B2").PasteSpecial(-4122)
$oldStatusSheet.Range("C3").Copy()
$newStatusSheet.Range("B3").PasteSpecial(-4122)
$oldStatusSheet.Range("C4").Copy()
$newStatusSheet.Range("B4").PasteSpecial(-4122)

# Match the new sheet's page setup (portrait, paper size 9 = A4).
$newStatusSheet.PageSetup.PaperSize = 9
$newStatusSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. Restore the recorded view selections for both "Status" sheets.
# ---------------------------------------------------------------------
$newStatusSheet.Range("N11").Select() | Out-Null
$oldStatusSheet.Range("B22").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Move the demoted sheet ("Status_1") to the end of the workbook,
#    after "API".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldStatusSheet.Move($null, $lastSheet)

# ---------------------------------------------------------------------
# 6. Add the new permission lookup entry on "Permissions" and leave it
#    as the active tab.
# ---------------------------------------------------------------------
$permissions = $wb.Worksheets.Item("Permissions")
$permissions.Range("A1").Copy()
$permissions.Range("F24").PasteSpecial(-4122)
$permissions.Range("F24").Value = "0,3"
$permissions.Range("F24").Select() | Out-Null
